$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $text) {
    $ws.Range("ZZ1").Formula = "=""" + $text + """"
    $ws.Range("ZZ1").Copy()
    $ws.Range($addr).PasteSpecial(-4163)
}

$ws.Range("D2").Value = "67.716.81"
$ws.Range("E2").Value = "  -1.61%  "
$ws.Range("D3").Value = "3.266.69"
$ws.Range("E3").Value = "  -0.75%  "
$ws.Range("E4").Value = "  -0.06%  "
Set-TextValue "D5" "580.60"
$ws.Range("E5").Value = "  -0.82%  "
Set-TextValue "D6" "184.83"
$ws.Range("E6").Value = "  +1.07%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("E8").Value = "  +0.59%  "
$ws.Range("E9").Value = "  -4.00%  "
$ws.Range("E10").Value = "  -1.35%  "
$ws.Range("E11").Value = "  -3.33%  "
$ws.Range("D12").Value = "3.828.02"
$ws.Range("E12").Value = "  -1.03%  "
$ws.Range("E13").Value = "  +0.49%  "
$ws.Range("E14").Value = "  -5.33%  "
$ws.Range("D15").Value = "67.772.29"
$ws.Range("E15").Value = "  -1.55%  "
$ws.Range("E16").Value = "  -2.35%  "
$ws.Range("D17").Value = "3.254.37"
$ws.Range("E17").Value = "  -3.18%  "
$ws.Range("E18").Value = "  -2.02%  "
Set-TextValue "D19" "13.45"
$ws.Range("E19").Value = "  -1.28%  "
Set-TextValue "D20" "398.68"
$ws.Range("E20").Value = "  +1.13%  "
$ws.Range("E21").Value = "  -2.29%  "
$ws.Range("E22").Value = "  +0.09%  "
Set-TextValue "D23" "70.78"
$ws.Range("E23").Value = "  -1.32%  "
Set-TextValue "D24" "0.509"
$ws.Range("E24").Value = "  -1.50%  "
$ws.Range("E25").Value = "  -3.50%  "
$ws.Range("E26").Value = "  -0.97%  "
$ws.Range("E27").Value = "  -1.85%  "
$ws.Range("E28").Value = "  +0.37%  "
$ws.Range("E29").Value = "  -1.83%  "
Set-TextValue "D30" "22.65"
$ws.Range("E30").Value = "  -1.96%  "
$ws.Range("E31").Value = "  -5.16%  "
$ws.Range("E32").Value = "  -3.11%  "
Set-TextValue "D33" "0.999"
$ws.Range("E33").Value = "  +0.05%  "
Set-TextValue "D34" "1.25"
$ws.Range("E34").Value = "  -4.23%  "
Set-TextValue "D35" "164.01"
$ws.Range("E35").Value = "  +0.21%  "
$ws.Range("E36").Value = "  -4.62%  "
$ws.Range("E37").Value = "  +0.06%  "
Set-TextValue "D38" "27.00"
$ws.Range("E38").Value = "  +2.37%  "
Set-TextValue "D39" "0.806"
$ws.Range("E39").Value = "  -3.56%  "
Set-TextValue "D40" "4.51"
$ws.Range("E40").Value = "  -2.19%  "
$ws.Range("D41").Value = "2.677.13"
$ws.Range("E41").Value = "  +2.36%  "
Set-TextValue "D42" "6.27"
$ws.Range("E42").Value = "  -4.76%  "
Set-TextValue "D43" "40.74"
$ws.Range("E43").Value = "  -2.01%  "
Set-TextValue "D44" "0.0678"
$ws.Range("E44").Value = "  -1.70%  "
Set-TextValue "D45" "2.43"
$ws.Range("E45").Value = "  -5.75%  "
Set-TextValue "D46" "335.35"
$ws.Range("E46").Value = "  -2.06%  "
Set-TextValue "D47" "24.56"
$ws.Range("E47").Value = "  -2.04%  "
$ws.Range("E48").Value = "  -3.35%  "
Set-TextValue "D49" "6.28"
$ws.Range("E49").Value = "  -0.48%  "
$ws.Range("E50").Value = "  -1.61%  "
Set-TextValue "D51" "0.969"
$ws.Range("E51").Value = "  -2.08%  "

$ws.Range("ZZ1").Clear()

